$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - F column "想去人数" updates
$wsExhibit = $wb.Worksheets.Item("展览")
$exhibitUpdates = @{
    2  = 1146
    3  = 1960
    4  = 618
    5  = 1269
    7  = 47
    8  = 139
    9  = 343
    10 = 125
    12 = 850
    13 = 260
    14 = 134
    15 = 33
    16 = 117
    18 = 251
    20 = 81
    21 = 672
    24 = 915
    25 = 370
    28 = 308
    31 = 428
}
foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Sheet "全部类型" (sheet4) - F column "想去人数" updates
$wsAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    3  = 1146
    4  = 1960
    5  = 618
    6  = 1269
    9  = 47
    10 = 139
    11 = 343
    12 = 125
    14 = 850
    15 = 260
    16 = 134
    18 = 33
    20 = 117
    25 = 251
    27 = 81
    28 = 672
    31 = 915
    32 = 370
    37 = 308
    43 = 428
}
foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
